# Auto-generated edit script applying the Raiden_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 1250068.5
$ws.Range("J5").Value = 9.199999999999999
$ws.Range("L5").Value = 9.199999999999999
$ws.Range("N5").Value = -239.2
$ws.Range("H17").Value = 1429.5
$ws.Range("J17").Value = 1429.5
$ws.Range("L17").Value = 4288.5
$ws.Range("N17").Value = -4624.5
$ws.Range("H51").Value = 9736.833000000001
$ws.Range("H53").Value = 750.93335
$ws.Range("J53").Value = 1630
$ws.Range("L53").Value = 1630
$ws.Range("N53").Value = -2904
$ws.Range("H62").Value = 4149.6665
$ws.Range("I62").Value = 4030.6667
$ws.Range("J62").Value = 4268.6665
$ws.Range("K62").Value = 4030.6667
$ws.Range("L62").Value = 4268.6665
$ws.Range("M62").Value = -3406.6667
$ws.Range("N62").Value = -5516.6665
$ws.Range("H65").Value = 4149.6665
$ws.Range("I65").Value = 4030.6667
$ws.Range("J65").Value = 4268.6665
$ws.Range("K65").Value = 20153.3335
$ws.Range("L65").Value = 21343.3325
$ws.Range("M65").Value = -17033.3335
$ws.Range("N65").Value = -27583.3325
$ws.Range("H88").Value = 669333
$ws.Range("I88").Value = 1999999
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 1999999
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -1999593
$ws.Range("N88").Value = -4812
$ws.Range("H91").Value = 669333
$ws.Range("I91").Value = 1999999
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 1999999
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = -1998595
$ws.Range("N91").Value = -6808
$ws.Range("H108").Value = 95227
$ws.Range("J108").Value = 95227
$ws.Range("L108").Value = 95227
$ws.Range("N108").Value = -102907
$ws.Range("H113").Value = 10838.857
$ws.Range("J113").Value = 14333
$ws.Range("L113").Value = 14333
$ws.Range("N113").Value = -20841
$ws.Range("H125").Value = 2341.4
$ws.Range("I125").Value = 2682
$ws.Range("J125").Value = 979
$ws.Range("K125").Value = 24138
$ws.Range("L125").Value = 8811
$ws.Range("M125").Value = -21678
$ws.Range("N125").Value = -13731
$ws.Range("H132").Value = 1216.814
$ws.Range("I132").Value = 1216.814
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3650.442
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1120.442
$ws.Range("N132").Value = $null
$ws.Range("H137").Value = 3209.1428
$ws.Range("I137").Value = 1597
$ws.Range("K137").Value = 4791
$ws.Range("M137").Value = -2241
$ws.Range("H138").Value = 2498.9604
$ws.Range("I138").Value = 5070
$ws.Range("J138").Value = 2196.4854
$ws.Range("K138").Value = 15210
$ws.Range("L138").Value = 6589.456200000001
$ws.Range("M138").Value = -10070
$ws.Range("N138").Value = -16869.4562

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 584.4
$ws.Range("I4").Value = 584.4
$ws.Range("K4").Value = 584.4
$ws.Range("M4").Value = -468.4
$ws.Range("H5").Value = 575
$ws.Range("I5").Value = 150
$ws.Range("K5").Value = 150
$ws.Range("M5").Value = -38
$ws.Range("H32").Value = 3405.5667
$ws.Range("I32").Value = 3178.1724
$ws.Range("K32").Value = 3178.1724
$ws.Range("M32").Value = -2891.1724
$ws.Range("H61").Value = 4091
$ws.Range("I61").Value = 4198
$ws.Range("K61").Value = 4198
$ws.Range("M61").Value = -3986
$ws.Range("H74").Value = 2462.2856
$ws.Range("I74").Value = 2462.2856
$ws.Range("K74").Value = 2462.2856
$ws.Range("M74").Value = -1588.2856
$ws.Range("H77").Value = 2462.2856
$ws.Range("I77").Value = 2462.2856
$ws.Range("K77").Value = 12311.428
$ws.Range("M77").Value = -7943.428
$ws.Range("H97").Value = 1455.3846
$ws.Range("I97").Value = 557.5714
$ws.Range("J97").Value = 2502.8333
$ws.Range("K97").Value = 557.5714
$ws.Range("L97").Value = 2502.8333
$ws.Range("M97").Value = -61.57140000000004
$ws.Range("N97").Value = -3494.8333
$ws.Range("H102").Value = 2399.6667
$ws.Range("I102").Value = 2399
$ws.Range("K102").Value = 2399
$ws.Range("M102").Value = -777
$ws.Range("H108").Value = 88821
$ws.Range("J108").Value = 88821
$ws.Range("L108").Value = 88821
$ws.Range("N108").Value = -96501
$ws.Range("H110").Value = 2106.6667
$ws.Range("I110").Value = 2106.6667
$ws.Range("K110").Value = 2106.6667
$ws.Range("M110").Value = -61.66670000000022
$ws.Range("H122").Value = 2243.8333
$ws.Range("J122").Value = 1549
$ws.Range("L122").Value = 4647
$ws.Range("N122").Value = -9547
$ws.Range("H132").Value = 3499.5
$ws.Range("I132").Value = 3998
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 11994
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -9464
$ws.Range("N132").Value = -15059.9999
$ws.Range("H136").Value = 4091
$ws.Range("I136").Value = 4198
$ws.Range("K136").Value = 12594
$ws.Range("M136").Value = -10044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 575
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -35
$ws.Range("H47").Value = 300000
$ws.Range("J47").Value = 300000
$ws.Range("L47").Value = 300000
$ws.Range("N47").Value = -301040
$ws.Range("H94").Value = 2382.6365
$ws.Range("I94").Value = 1091.4
$ws.Range("K94").Value = 1091.4
$ws.Range("M94").Value = -640.4000000000001
$ws.Range("H107").Value = 749.7
$ws.Range("I107").Value = 749.625
$ws.Range("J107").Value = 750
$ws.Range("K107").Value = 749.625
$ws.Range("L107").Value = 750
$ws.Range("M107").Value = 1170.375
$ws.Range("N107").Value = -4590
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 492.66666
$ws.Range("I22").Value = 500.66666
$ws.Range("J22").Value = 476.66666
$ws.Range("K22").Value = 500.66666
$ws.Range("L22").Value = 476.66666
$ws.Range("M22").Value = -150.66666
$ws.Range("N22").Value = -1176.66666
$ws.Range("H31").Value = 3664.7778
$ws.Range("I31").Value = 2207.6
$ws.Range("K31").Value = 2207.6
$ws.Range("M31").Value = -1912.6
$ws.Range("H34").Value = 3664.7778
$ws.Range("I34").Value = 2207.6
$ws.Range("K34").Value = 2207.6
$ws.Range("M34").Value = -2005.6
$ws.Range("H58").Value = 3289.2
$ws.Range("I58").Value = 3374.375
$ws.Range("K58").Value = 3374.375
$ws.Range("M58").Value = -3171.375
$ws.Range("H74").Value = 64974.75
$ws.Range("J74").Value = 64974.75
$ws.Range("L74").Value = 64974.75
$ws.Range("N74").Value = -66722.75
$ws.Range("H77").Value = 64974.75
$ws.Range("J77").Value = 64974.75
$ws.Range("L77").Value = 194924.25
$ws.Range("N77").Value = -203660.25
$ws.Range("H81").Value = 62000
$ws.Range("J81").Value = 62000
$ws.Range("L81").Value = 62000
$ws.Range("N81").Value = -63996
$ws.Range("H84").Value = 62000
$ws.Range("J84").Value = 62000
$ws.Range("L84").Value = 186000
$ws.Range("N84").Value = -195984
$ws.Range("H99").Value = 1670.3334
$ws.Range("J99").Value = 1499
$ws.Range("L99").Value = 1499
$ws.Range("N99").Value = -4495
$ws.Range("H105").Value = 1913.1666
$ws.Range("I105").Value = 1474.2858
$ws.Range("J105").Value = 3449.25
$ws.Range("K105").Value = 1474.2858
$ws.Range("L105").Value = 3449.25
$ws.Range("M105").Value = 272.7141999999999
$ws.Range("N105").Value = -6943.25
$ws.Range("H107").Value = 1135.6842
$ws.Range("I107").Value = 1244.3636
$ws.Range("J107").Value = 986.25
$ws.Range("K107").Value = 1244.3636
$ws.Range("L107").Value = 986.25
$ws.Range("M107").Value = 675.6364000000001
$ws.Range("N107").Value = -4826.25
$ws.Range("H122").Value = 1673.75
$ws.Range("I122").Value = 1627
$ws.Range("K122").Value = 4881
$ws.Range("M122").Value = -2431
$ws.Range("H126").Value = 1670.3334
$ws.Range("J126").Value = 1499
$ws.Range("L126").Value = 4497
$ws.Range("N126").Value = -9437
$ws.Range("H132").Value = 2017.3684
$ws.Range("I132").Value = 2017.3684
$ws.Range("K132").Value = 6052.1052
$ws.Range("M132").Value = -3522.1052
$ws.Range("H134").Value = 1824.7222
$ws.Range("I134").Value = 1998.0769
$ws.Range("K134").Value = 5994.2307
$ws.Range("M134").Value = -3459.2307
$ws.Range("H136").Value = 3289.2
$ws.Range("I136").Value = 3374.375
$ws.Range("K136").Value = 10123.125
$ws.Range("M136").Value = -7573.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 333377.72
$ws.Range("I2").Value = 357187.16
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 2143122.96
$ws.Range("L2").Value = 276
$ws.Range("M2").Value = -2143009.96
$ws.Range("N2").Value = -502
$ws.Range("I4").Value = 4198212.5
$ws.Range("J4").Value = 66667056
$ws.Range("K4").Value = 12594637.5
$ws.Range("L4").Value = 200001168
$ws.Range("M4").Value = -12594525.5
$ws.Range("N4").Value = -200001392
$ws.Range("H6").Value = 115
$ws.Range("I6").Value = 115
$ws.Range("K6").Value = 345
$ws.Range("M6").Value = -232
$ws.Range("H12").Value = 346.1905
$ws.Range("I12").Value = 428.8
$ws.Range("J12").Value = 271.0909
$ws.Range("K12").Value = 1286.4
$ws.Range("L12").Value = 813.2727
$ws.Range("M12").Value = -1113.4
$ws.Range("N12").Value = -1159.2727
$ws.Range("H80").Value = 25364.883
$ws.Range("I80").Value = 13699.143
$ws.Range("J80").Value = 33530.9
$ws.Range("K80").Value = 41097.429
$ws.Range("L80").Value = 100592.7
$ws.Range("M80").Value = -40161.429
$ws.Range("N80").Value = -102464.7
$ws.Range("H83").Value = 25364.883
$ws.Range("I83").Value = 13699.143
$ws.Range("J83").Value = 33530.9
$ws.Range("K83").Value = 123292.287
$ws.Range("L83").Value = 301778.1
$ws.Range("M83").Value = -118612.287
$ws.Range("N83").Value = -311138.1
$ws.Range("H113").Value = 812.7778
$ws.Range("I113").Value = 638.75
$ws.Range("J113").Value = 952
$ws.Range("K113").Value = 1916.25
$ws.Range("L113").Value = 2856
$ws.Range("M113").Value = 253.75
$ws.Range("N113").Value = -7196
$ws.Range("H132").Value = 874.625
$ws.Range("I132").Value = 666.1667
$ws.Range("K132").Value = 5995.5003
$ws.Range("M132").Value = -3465.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 9126.75
$ws.Range("I9").Value = 500
$ws.Range("J9").Value = 12002.333
$ws.Range("K9").Value = 500
$ws.Range("L9").Value = 12002.333
$ws.Range("M9").Value = -330
$ws.Range("N9").Value = -12342.333
$ws.Range("H29").Value = 1249.5
$ws.Range("I29").Value = 1249.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1249.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -959.5
$ws.Range("N29").Value = $null
$ws.Range("H107").Value = 590.9474
$ws.Range("I107").Value = 267.36365
$ws.Range("J107").Value = 1035.875
$ws.Range("K107").Value = 267.36365
$ws.Range("L107").Value = 1035.875
$ws.Range("M107").Value = 1652.63635
$ws.Range("N107").Value = -4875.875
$ws.Range("H126").Value = 3059.8235
$ws.Range("J126").Value = 3047.2222
$ws.Range("L126").Value = 9141.6666
$ws.Range("N126").Value = -14081.6666
$ws.Range("H132").Value = 2414.75
$ws.Range("I132").Value = 2420.3333
$ws.Range("K132").Value = 7260.999899999999
$ws.Range("M132").Value = -4730.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H16").Value = 315.2353
$ws.Range("I16").Value = 330.66666
$ws.Range("K16").Value = 330.66666
$ws.Range("M16").Value = -160.66666
$ws.Range("H42").Value = 40000.715
$ws.Range("I42").Value = 36667.5
$ws.Range("K42").Value = 36667.5
$ws.Range("M42").Value = -36104.5
$ws.Range("H49").Value = 40000.715
$ws.Range("I49").Value = 36667.5
$ws.Range("K49").Value = 36667.5
$ws.Range("M49").Value = -36520.5
$ws.Range("H61").Value = 987.9231
$ws.Range("I61").Value = 978.8333
$ws.Range("J61").Value = 1097
$ws.Range("K61").Value = 978.8333
$ws.Range("L61").Value = 1097
$ws.Range("M61").Value = -776.8333
$ws.Range("N61").Value = -1501
$ws.Range("H68").Value = 2475.65
$ws.Range("I68").Value = 2595.3076
$ws.Range("J68").Value = 2253.4285
$ws.Range("K68").Value = 2595.3076
$ws.Range("L68").Value = 2253.4285
$ws.Range("M68").Value = -1846.3076
$ws.Range("N68").Value = -3751.4285
$ws.Range("H71").Value = 2475.65
$ws.Range("I71").Value = 2595.3076
$ws.Range("J71").Value = 2253.4285
$ws.Range("K71").Value = 12976.538
$ws.Range("L71").Value = 11267.1425
$ws.Range("M71").Value = -9232.538
$ws.Range("N71").Value = -18755.1425
$ws.Range("H93").Value = 943.3
$ws.Range("I93").Value = 945.8333
$ws.Range("J93").Value = 939.5
$ws.Range("K93").Value = 945.8333
$ws.Range("L93").Value = 939.5
$ws.Range("M93").Value = 302.1667
$ws.Range("N93").Value = -3435.5
$ws.Range("H113").Value = 987.9231
$ws.Range("I113").Value = 978.8333
$ws.Range("J113").Value = 1097
$ws.Range("K113").Value = 978.8333
$ws.Range("L113").Value = 1097
$ws.Range("M113").Value = 1191.1667
$ws.Range("N113").Value = -5437
$ws.Range("H122").Value = 3685.879
$ws.Range("I122").Value = 3712.7827
$ws.Range("J122").Value = 3624
$ws.Range("K122").Value = 11138.3481
$ws.Range("L122").Value = 10872
$ws.Range("M122").Value = -8688.348100000001
$ws.Range("N122").Value = -15772
$ws.Range("H136").Value = 1471.8462
$ws.Range("I136").Value = 1372.9714
$ws.Range("J136").Value = 2337
$ws.Range("K136").Value = 4118.914199999999
$ws.Range("L136").Value = 7011
$ws.Range("M136").Value = -1568.914199999999
$ws.Range("N136").Value = -12111
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 31459.143
$ws.Range("J43").Value = 35612
$ws.Range("L43").Value = 35612
$ws.Range("N43").Value = -35910
$ws.Range("H68").Value = 13758.5
$ws.Range("I68").Value = 11246
$ws.Range("J68").Value = 16271
$ws.Range("K68").Value = 11246
$ws.Range("L68").Value = 16271
$ws.Range("M68").Value = -10435
$ws.Range("N68").Value = -17893
$ws.Range("H69").Value = 44840
$ws.Range("J69").Value = 44840
$ws.Range("L69").Value = 44840
$ws.Range("N69").Value = -46338
$ws.Range("H70").Value = 18045.182
$ws.Range("I70").Value = 15570.571
$ws.Range("J70").Value = 19200
$ws.Range("K70").Value = 15570.571
$ws.Range("L70").Value = 19200
$ws.Range("M70").Value = -15255.571
$ws.Range("N70").Value = -19830
$ws.Range("H71").Value = 13758.5
$ws.Range("I71").Value = 11246
$ws.Range("J71").Value = 16271
$ws.Range("K71").Value = 33738
$ws.Range("L71").Value = 48813
$ws.Range("M71").Value = -29682
$ws.Range("N71").Value = -56925
$ws.Range("H72").Value = 44840
$ws.Range("J72").Value = 44840
$ws.Range("L72").Value = 134520
$ws.Range("N72").Value = -142008
$ws.Range("H73").Value = 18045.182
$ws.Range("I73").Value = 15570.571
$ws.Range("J73").Value = 19200
$ws.Range("K73").Value = 15570.571
$ws.Range("L73").Value = 19200
$ws.Range("M73").Value = -14478.571
$ws.Range("N73").Value = -21384
$ws.Range("H96").Value = 2426623.5
$ws.Range("I96").Value = 3032780
$ws.Range("J96").Value = 1998
$ws.Range("K96").Value = 3032780
$ws.Range("L96").Value = 1998
$ws.Range("M96").Value = -3031407
$ws.Range("N96").Value = -4744
$ws.Range("H99").Value = 56665.332
$ws.Range("I99").Value = 56665.332
$ws.Range("K99").Value = 56665.332
$ws.Range("M99").Value = -53670.332
$ws.Range("H100").Value = 1673.5
$ws.Range("I100").Value = 2098.3333
$ws.Range("J100").Value = 399
$ws.Range("K100").Value = 4196.6666
$ws.Range("L100").Value = 798
$ws.Range("M100").Value = -3655.6666
$ws.Range("N100").Value = -1880
$ws.Range("H113").Value = 408.46155
$ws.Range("I113").Value = 341.3
$ws.Range("K113").Value = 1023.9
$ws.Range("M113").Value = 1146.1
$ws.Range("H122").Value = 6581.619
$ws.Range("I122").Value = 8519
$ws.Range("J122").Value = 4450.5
$ws.Range("K122").Value = 25557
$ws.Range("L122").Value = 13351.5
$ws.Range("M122").Value = -23107
$ws.Range("N122").Value = -18251.5
$ws.Range("H126").Value = 4068.35
$ws.Range("I126").Value = 4172.8125
$ws.Range("J126").Value = 3650.5
$ws.Range("K126").Value = 12518.4375
$ws.Range("L126").Value = 10951.5
$ws.Range("M126").Value = -10048.4375
$ws.Range("N126").Value = -15891.5
$ws.Range("H132").Value = 1962
$ws.Range("I132").Value = 1917.6923
$ws.Range("K132").Value = 5753.0769
$ws.Range("M132").Value = -3223.0769
$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360

